$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same rows and need the
# same "想去人数" (F column) updates.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 2186
    $ws.Range("F3").Value = 626
    $ws.Range("F4").Value = 1567
    $ws.Range("F5").Value = 7337
    $ws.Range("F7").Value = 178
}
